$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1. Title: "Pi2Go Simulator Programming: " -> "Virtual Pi2Go Programming: "
# ---------------------------------------------------------------
$d.Content.Find.Execute("Pi2Go Simulator Programming: ", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Virtual Pi2Go Programming: ", 2) | Out-Null

# ---------------------------------------------------------------
# 2. "It checks to see if any rules are applicable" -> insert "BDI "
#    before "rules " (2nd occurrence of the lead-in phrase).
# ---------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("It checks to see if any ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Start = $rng.End
$rng.End = $d.Content.End
$rng.Find.Execute("It checks to see if any ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.InsertAfter("BDI ")

# ---------------------------------------------------------------
# 3. "What is a rule?  A rule is ..." -> "What is a BDI rule?  A BDI rule is ..."
# ---------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("What is a rule?  A rule is", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$base = $rng.Start
$r2 = $d.Range($base + 19, $base + 19)
$r2.InsertBefore("BDI ")
$r1 = $d.Range($base + 10, $base + 10)
$r1.InsertBefore("BDI ")

# ---------------------------------------------------------------
# 4. Remove one empty paragraph right after "Rules are always applicable..."
# ---------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Rules are always applicable unless they have a condition.  We will talk about conditions in the next section.", `
                   $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$targetIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $rng.Start -and $rng.Start -lt $p.Range.End) {
        $targetIdx = $i
        break
    }
}
$nextP = $d.Paragraphs.Item($targetIdx + 1)
$nextP.Range.Delete()
